# "Formal methods started, added to SA process."
#
# 1) The cached "today" value of the datetimeFigureOut field (shown on the
#    slide master + every slide layout's Date placeholder) moved on from
#    12/6/2019 to 12/9/2019 - the deck was simply reopened a few days later.
# 2) The V-model safety-assessment diagram on slide 2 was updated: each
#    box's activity label now calls out the specific techniques (FTA, CMA,
#    FHA, CCA, PSSA, PASA, ASA, SSA, FMEA) instead of the generic assessment
#    name, reflecting that formal methods have been folded into the process.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "date last opened" field everywhere it is cached:
#    the slide master and all eleven slide layouts.
# ---------------------------------------------------------------------------
function Update-CachedDate($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if (-not $sh.HasTextFrame) { continue }
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "12/6/2019") {
            $tr.Text = "12/9/2019"
        }
    }
}

$master = $p.SlideMaster
Update-CachedDate $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    Update-CachedDate $master.CustomLayouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 2 - the V-model diagram's rounded-rectangle activity boxes.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)

# "Preliminary Module Safety Assessment" -> "System FTA, CMA"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "System FTA, CMA"

# "Preliminary System Safety Assessment" -> three stacked lines:
#   System FHA
#   CCA
#   PSSA
$s.Shapes.Item(4).TextFrame.TextRange.Text = "System FHA`rCCA`rPSSA"

# "Functional Hazard Assessment" -> three stacked lines, the first one built
# from two runs ("Aircraft " + "FHA"), followed by "CCA" and "PASA"
$tr5 = $s.Shapes.Item(5).TextFrame.TextRange
$tr5.Text = "Aircraft "
[void]$tr5.InsertAfter("FHA")
[void]$tr5.InsertAfter("`rCCA`rPASA")

# "Safety Certification" -> "Aircraft CCA, ASA"
$s.Shapes.Item(6).TextFrame.TextRange.Text = "Aircraft CCA, ASA"

# "System Safety Assessment" -> "System SSA, CCA, FMEA"
$s.Shapes.Item(7).TextFrame.TextRange.Text = "System SSA, CCA, FMEA"

# "Module Safety Assessment" -> "System FTA, CMA, FMEA"
$s.Shapes.Item(8).TextFrame.TextRange.Text = "System FTA, CMA, FMEA"
